$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 4.92663022216654
$ws.Cells.Item(2, 4).Value = 2.706030255740866
$ws.Cells.Item(2, 5).Value = 40.45870857033537
$ws.Cells.Item(2, 6).Value = 20.63069856659229
$ws.Cells.Item(2, 7).Value = 25.07669775851724
$ws.Cells.Item(2, 8).Value = 11.75330198768288
$ws.Cells.Item(2, 13).Value = 57.86521735402618
$ws.Cells.Item(3, 3).Value = 4.75020579027044
$ws.Cells.Item(3, 4).Value = 2.775981400822221
$ws.Cells.Item(3, 5).Value = 37.7447454554132
$ws.Cells.Item(3, 6).Value = 20.63698421196751
$ws.Cells.Item(3, 7).Value = 24.63225408136648
$ws.Cells.Item(3, 8).Value = 11.95258782237333
$ws.Cells.Item(3, 13).Value = 54.32657143557415
$ws.Cells.Item(4, 3).Value = 4.640721872203416
$ws.Cells.Item(4, 4).Value = 2.882741088580405
$ws.Cells.Item(4, 5).Value = 35.97835745701322
$ws.Cells.Item(4, 6).Value = 20.68308875040475
$ws.Cells.Item(4, 7).Value = 24.43684844615833
$ws.Cells.Item(4, 8).Value = 12.08886615363643
$ws.Cells.Item(4, 13).Value = 52.03338129302649
$ws.Cells.Item(5, 3).Value = 4.595888422962744
$ws.Cells.Item(5, 4).Value = 2.927040347804819
$ws.Cells.Item(5, 5).Value = 35.2332020867936
$ws.Cells.Item(5, 6).Value = 20.71204526999227
$ws.Cells.Item(5, 7).Value = 24.37629344396692
$ws.Cells.Item(5, 8).Value = 12.14774615468647
$ws.Cells.Item(5, 13).Value = 51.06868264551071
$ws.Cells.Item(6, 3).Value = 4.588432967093631
$ws.Cells.Item(6, 4).Value = 2.93444452074627
$ws.Cells.Item(6, 5).Value = 35.10793370229888
$ws.Cells.Item(6, 6).Value = 20.71745366985141
$ws.Cells.Item(6, 7).Value = 24.36737443592509
$ws.Cells.Item(6, 8).Value = 12.15772086564211
$ws.Cells.Item(6, 13).Value = 50.90667464537157
$ws.Cells.Item(7, 3).Value = 4.640118013298389
$ws.Cells.Item(7, 4).Value = 2.883335292988268
$ws.Cells.Item(7, 5).Value = 35.96841087565716
$ws.Cells.Item(7, 6).Value = 20.68343872255418
$ws.Cells.Item(7, 7).Value = 24.43595524023551
$ws.Cells.Item(7, 8).Value = 12.08964688093299
$ws.Cells.Item(7, 13).Value = 52.0204931698381
$ws.Cells.Item(8, 3).Value = 4.866081757249129
$ws.Cells.Item(8, 4).Value = 2.664736848931564
$ws.Cells.Item(8, 5).Value = 39.54357401156243
$ws.Cells.Item(8, 6).Value = 20.62383684784165
$ws.Cells.Item(8, 7).Value = 24.90718999635376
$ws.Cells.Item(8, 8).Value = 11.81903217340908
$ws.Cells.Item(8, 13).Value = 56.67011963173081
$ws.Cells.Item(9, 3).Value = 5.297091858814765
$ws.Cells.Item(9, 4).Value = 3.038493853217551
$ws.Cells.Item(9, 5).Value = 45.76969450842201
$ws.Cells.Item(9, 6).Value = 20.86007096918814
$ws.Cells.Item(9, 7).Value = 26.45342108163448
$ws.Cells.Item(9, 8).Value = 11.40612644023092
$ws.Cells.Item(9, 13).Value = 64.83201972146196
$ws.Cells.Item(10, 3).Value = 5.602939013525972
$ws.Cells.Item(10, 4).Value = 3.310371983115453
$ws.Cells.Item(10, 5).Value = 49.87786705760199
$ws.Cells.Item(10, 6).Value = 21.26918004761591
$ws.Cells.Item(10, 7).Value = 27.96451690518017
$ws.Cells.Item(10, 8).Value = 11.18532037229778
$ws.Cells.Item(10, 13).Value = 70.2492148944378
$ws.Cells.Item(11, 3).Value = 5.739132207954444
$ws.Cells.Item(11, 4).Value = 3.431539541104398
$ws.Cells.Item(11, 5).Value = 51.64827665769101
$ws.Cells.Item(11, 6).Value = 21.50944500877759
$ws.Cells.Item(11, 7).Value = 28.72917904519621
$ws.Cells.Item(11, 8).Value = 11.10530132852388
$ws.Cells.Item(11, 13).Value = 72.58943003491193
$ws.Cells.Item(12, 3).Value = 5.79024120338491
$ws.Cells.Item(12, 4).Value = 3.477086280605161
$ws.Cells.Item(12, 5).Value = 52.30475915911518
$ws.Cells.Item(12, 6).Value = 21.60836661078242
$ws.Cells.Item(12, 7).Value = 29.02934808015367
$ws.Cells.Item(12, 8).Value = 11.07815567630058
$ws.Cells.Item(12, 13).Value = 73.4579325273275
$ws.Cells.Item(13, 3).Value = 5.77925522595607
$ws.Cells.Item(13, 4).Value = 3.467291549873048
$ws.Cells.Item(13, 5).Value = 52.16399002429173
$ws.Cells.Item(13, 6).Value = 21.58670733174029
$ws.Cells.Item(13, 7).Value = 28.96423811760813
$ws.Cells.Item(13, 8).Value = 11.08385806948573
$ws.Cells.Item(13, 13).Value = 73.27166875199566
$ws.Cells.Item(14, 3).Value = 5.74334650223881
$ws.Cells.Item(14, 4).Value = 3.435293303567023
$ws.Cells.Item(14, 5).Value = 51.70256383002291
$ws.Cells.Item(14, 6).Value = 21.51742362863297
$ws.Cells.Item(14, 7).Value = 28.75366461633508
$ws.Cells.Item(14, 8).Value = 11.10300351563863
$ws.Cells.Item(14, 13).Value = 72.66123570762353
$ws.Cells.Item(15, 3).Value = 5.721289749658454
$ws.Cells.Item(15, 4).Value = 3.415650342826378
$ws.Cells.Item(15, 5).Value = 51.41811907229584
$ws.Cells.Item(15, 6).Value = 21.47602249071823
$ws.Cells.Item(15, 7).Value = 28.62604781661805
$ws.Cells.Item(15, 8).Value = 11.11514824066037
$ws.Cells.Item(15, 13).Value = 72.28502976236319
$ws.Cells.Item(16, 3).Value = 5.593975625664173
$ws.Cells.Item(16, 4).Value = 3.302404910598212
$ws.Cells.Item(16, 5).Value = 49.76020478707021
$ws.Cells.Item(16, 6).Value = 21.25458240661795
$ws.Cells.Item(16, 7).Value = 27.91605771326773
$ws.Cells.Item(16, 8).Value = 11.19097935317579
$ws.Cells.Item(16, 13).Value = 70.09379212449794
$ws.Cells.Item(17, 3).Value = 5.515088271466781
$ws.Cells.Item(17, 4).Value = 3.232306153207265
$ws.Cells.Item(17, 5).Value = 48.71807467179959
$ws.Cells.Item(17, 6).Value = 21.13273631668204
$ws.Cells.Item(17, 7).Value = 27.49995508689101
$ws.Cells.Item(17, 8).Value = 11.2428841432117
$ws.Cells.Item(17, 13).Value = 68.71785881290143
$ws.Cells.Item(18, 3).Value = 5.469439908510362
$ws.Cells.Item(18, 4).Value = 3.191749292772599
$ws.Cells.Item(18, 5).Value = 48.10939186805792
$ws.Cells.Item(18, 6).Value = 21.06774679043274
$ws.Cells.Item(18, 7).Value = 27.2679314566601
$ws.Cells.Item(18, 8).Value = 11.27464748165998
$ws.Cells.Item(18, 13).Value = 67.91476662660979
$ws.Cells.Item(19, 3).Value = 5.453938395992205
$ws.Cells.Item(19, 4).Value = 3.177975859730464
$ws.Cells.Item(19, 5).Value = 47.90170239907148
$ws.Cells.Item(19, 6).Value = 21.04661149648404
$ws.Cells.Item(19, 7).Value = 27.1906414706729
$ws.Cells.Item(19, 8).Value = 11.28572335226386
$ws.Cells.Item(19, 13).Value = 67.64084099078268
$ws.Cells.Item(20, 3).Value = 5.523514691258606
$ws.Cells.Item(20, 4).Value = 3.239792737069351
$ws.Cells.Item(20, 5).Value = 48.82996978100088
$ws.Cells.Item(20, 6).Value = 21.14517861547516
$ws.Cells.Item(20, 7).Value = 27.54349703904396
$ws.Cells.Item(20, 8).Value = 11.23715964808365
$ws.Cells.Item(20, 13).Value = 68.86553825618023
$ws.Cells.Item(21, 3).Value = 5.753906673425405
$ws.Cells.Item(21, 4).Value = 3.444700902800935
$ws.Cells.Item(21, 5).Value = 51.83847224999617
$ws.Cells.Item(21, 6).Value = 21.53755762056668
$ws.Cells.Item(21, 7).Value = 28.81523155187693
$ws.Cells.Item(21, 8).Value = 11.09729260996774
$ws.Cells.Item(21, 13).Value = 72.84101328552929
$ws.Cells.Item(22, 3).Value = 5.901758879272341
$ws.Cells.Item(22, 4).Value = 3.576669624231359
$ws.Cells.Item(22, 5).Value = 53.72359598818063
$ws.Cells.Item(22, 6).Value = 21.84029416659346
$ws.Cells.Item(22, 7).Value = 29.70794870866591
$ws.Cells.Item(22, 8).Value = 11.02438190045184
$ws.Cells.Item(22, 13).Value = 75.33622031797186
$ws.Cells.Item(23, 3).Value = 5.823108791016317
$ws.Cells.Item(23, 4).Value = 3.506405465703755
$ws.Cells.Item(23, 5).Value = 52.72481697347003
$ws.Cells.Item(23, 6).Value = 21.67444968094836
$ws.Cells.Item(23, 7).Value = 29.2260349566339
$ws.Cells.Item(23, 8).Value = 11.06152901608675
$ws.Cells.Item(23, 13).Value = 74.0138459480233
$ws.Cells.Item(24, 3).Value = 5.51970602239843
$ws.Cells.Item(24, 4).Value = 3.236408852178645
$ws.Cells.Item(24, 5).Value = 48.77941175952597
$ws.Cells.Item(24, 6).Value = 21.13953769922808
$ws.Cells.Item(24, 7).Value = 27.52378926061497
$ws.Cells.Item(24, 8).Value = 11.23974172038679
$ws.Cells.Item(24, 13).Value = 68.79880992397341
$ws.Cells.Item(25, 3).Value = 5.182168030038254
$ws.Cells.Item(25, 4).Value = 2.935991824725725
$ws.Cells.Item(25, 5).Value = 44.16812726076083
$ws.Cells.Item(25, 6).Value = 20.756035967963
$ws.Cells.Item(25, 7).Value = 25.96826590553401
$ws.Cells.Item(25, 8).Value = 11.50414683576907
$ws.Cells.Item(25, 13).Value = 62.72616813052861
